$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the worker/period detail table (B16:G33) so it is ordered by
# Periodo Mora ascending (1810, 1811, 1812, 1901, 1902, 1903), and within
# each period by worker (FABIO 73009947, MARYSEL 45478050, OSCAR 73572972).
# Also corrects FABIO's 1903 "Valor Mora" from 40000 to 48000.
# The rest of the underlying (worker, period) -> (Valor Mora, Salario
# Basico) facts are unchanged; only the row order changes.

$data = @(
    @("CC", "73009947", "FABIO ANDRES ARROYO BELTRAN", "1810", 48000, 1200000),
    @("CC", "45478050", "MARYSEL CAÑAS PALACIO",        "1810", 32000, 800000),
    @("CC", "73572972", "OSCAR LUIS MENDIETA ESTARITA",  "1810", 60000, 1500000),
    @("CC", "73009947", "FABIO ANDRES ARROYO BELTRAN", "1811", 48000, 1200000),
    @("CC", "45478050", "MARYSEL CAÑAS PALACIO",        "1811", 32000, 800000),
    @("CC", "73572972", "OSCAR LUIS MENDIETA ESTARITA",  "1811", 60000, 1500000),
    @("CC", "73009947", "FABIO ANDRES ARROYO BELTRAN", "1812", 48000, 1200000),
    @("CC", "45478050", "MARYSEL CAÑAS PALACIO",        "1812", 32000, 800000),
    @("CC", "73572972", "OSCAR LUIS MENDIETA ESTARITA",  "1812", 60000, 1500000),
    @("CC", "73009947", "FABIO ANDRES ARROYO BELTRAN", "1901", 48000, 1200000),
    @("CC", "45478050", "MARYSEL CAÑAS PALACIO",        "1901", 32000, 800000),
    @("CC", "73572972", "OSCAR LUIS MENDIETA ESTARITA",  "1901", 60000, 1500000),
    @("CC", "73009947", "FABIO ANDRES ARROYO BELTRAN", "1902", 48000, 1200000),
    @("CC", "45478050", "MARYSEL CAÑAS PALACIO",        "1902", 32000, 800000),
    @("CC", "73572972", "OSCAR LUIS MENDIETA ESTARITA",  "1902", 60000, 1500000),
    @("CC", "73009947", "FABIO ANDRES ARROYO BELTRAN", "1903", 40000, 1200000),
    @("CC", "45478050", "MARYSEL CAÑAS PALACIO",        "1903", 26667, 800000),
    @("CC", "73572972", "OSCAR LUIS MENDIETA ESTARITA",  "1903", 50000, 1500000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
